$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.184.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.901.14'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.09%  '

$ws.Range("E4").Value = '  -0.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.95%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.697'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.72%  '

$ws.Range("E7").Value = '  -0.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.99'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.64%  '

$ws.Range("E9").Value = '  +3.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.92'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0753'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.77%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0983'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.18%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '13.03'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.176.08'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.737'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.25%  '

$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.96'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.55%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.876.60'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.43%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '35.166.88'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.31%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.77%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0834'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '242.45'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.57%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.96'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.81%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.50%  '

$ws.Range("E24").Value = '  -0.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.28'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.55%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.59'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.71%  '

$ws.Range("E30").Value = '  -0.92%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.128.51'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.55%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.08'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +16.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0608'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.33'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.61%  '

$ws.Range("B35").Value = 'TrustWalletToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.56'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +16.99%  '

$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.21'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.14%  '

$ws.Range("E37").Value = '  -0.29%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.853'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -10.70%  '

$ws.Range("E39").Value = '  -1.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '101.29'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +12.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.24'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.66%  '

$ws.Range("E42").Value = '  +1.85%  '

$ws.Range("E43").Value = '  +0.12%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0649'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.79%  '

$ws.Range("E45").Value = '  -0.45%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.315.09'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.30%  '

$ws.Range("E47").Value = '  +0.28%  '

$ws.Range("E48").Value = '  -1.74%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.84'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.88%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0740'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.63%  '
